$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet — "想去人数" (F column) count updates
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 54937
$wsExhibit.Range("F6").Value = 1172
$wsExhibit.Range("F10").Value = 1096
$wsExhibit.Range("F11").Value = 1359
$wsExhibit.Range("F14").Value = 215
$wsExhibit.Range("F16").Value = 54
$wsExhibit.Range("F21").Value = 5384
$wsExhibit.Range("F23").Value = 5273
$wsExhibit.Range("F24").Value = 9273
$wsExhibit.Range("F28").Value = 238
$wsExhibit.Range("F29").Value = 449
$wsExhibit.Range("F32").Value = 4270
$wsExhibit.Range("F33").Value = 283

# 演出 (Performance) sheet — "想去人数" (F column) count update
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F10").Value = 1149

# 全部类型 (All types) sheet — "想去人数" (F column) count updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 1172
$wsAll.Range("F13").Value = 1096
$wsAll.Range("F15").Value = 1359
$wsAll.Range("F18").Value = 215
$wsAll.Range("F21").Value = 54
$wsAll.Range("F26").Value = 5384
$wsAll.Range("F28").Value = 5273
$wsAll.Range("F29").Value = 9273
$wsAll.Range("F34").Value = 238
$wsAll.Range("F35").Value = 449
$wsAll.Range("F41").Value = 4270
$wsAll.Range("F47").Value = 283
